$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 24841.691
$ws.Range("I137").Value = 44447.957
$ws.Range("K137").Value = 133343.871
$ws.Range("M137").Value = -130793.871
$ws.Range("H138").Value = 1485.5858
$ws.Range("I138").Value = 802.0857
$ws.Range("J138").Value = 1859.375
$ws.Range("K138").Value = 2406.2571
$ws.Range("L138").Value = 5578.125
$ws.Range("M138").Value = 2733.7429
$ws.Range("N138").Value = -15858.125

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17212.834
$ws.Range("I2").Value = 1054.2
$ws.Range("J2").Value = 98006
$ws.Range("K2").Value = 1054.2
$ws.Range("L2").Value = 98006
$ws.Range("M2").Value = -941.2
$ws.Range("N2").Value = -98232
$ws.Range("H32").Value = 22809.5
$ws.Range("I32").Value = 5726.46
$ws.Range("J32").Value = 129578.5
$ws.Range("K32").Value = 5726.46
$ws.Range("L32").Value = 129578.5
$ws.Range("M32").Value = -5439.46
$ws.Range("N32").Value = -130152.5
$ws.Range("H45").Value = 1232.4
$ws.Range("I45").Value = 1106
$ws.Range("J45").Value = 1316.6666
$ws.Range("K45").Value = 1106
$ws.Range("L45").Value = 1316.6666
$ws.Range("M45").Value = -729
$ws.Range("N45").Value = -2070.6666
$ws.Range("H61").Value = 11905518
$ws.Range("I61").Value = 12500706
$ws.Range("J61").Value = 1750
$ws.Range("K61").Value = 12500706
$ws.Range("L61").Value = 1750
$ws.Range("M61").Value = -12500494
$ws.Range("N61").Value = -2174
$ws.Range("H102").Value = 1792.2307
$ws.Range("I102").Value = 1663.5454
$ws.Range("K102").Value = 1663.5454
$ws.Range("M102").Value = -41.54539999999997
$ws.Range("H116").Value = 17212.834
$ws.Range("I116").Value = 1054.2
$ws.Range("J116").Value = 98006
$ws.Range("K116").Value = 1054.2
$ws.Range("L116").Value = 98006
$ws.Range("M116").Value = 1239.8
$ws.Range("N116").Value = -102594
$ws.Range("H136").Value = 11905518
$ws.Range("I136").Value = 12500706
$ws.Range("J136").Value = 1750
$ws.Range("K136").Value = 37502118
$ws.Range("L136").Value = 5250
$ws.Range("M136").Value = -37499568
$ws.Range("N136").Value = -10350

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17212.834
$ws.Range("I3").Value = 1054.2
$ws.Range("J3").Value = 98006
$ws.Range("K3").Value = 1054.2
$ws.Range("L3").Value = 98006
$ws.Range("M3").Value = -940.2
$ws.Range("N3").Value = -98234

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28111.975
$ws.Range("I31").Value = 37019.965
$ws.Range("J31").Value = 4627.273
$ws.Range("K31").Value = 37019.965
$ws.Range("L31").Value = 4627.273
$ws.Range("M31").Value = -36724.965
$ws.Range("N31").Value = -5217.273
$ws.Range("H34").Value = 28111.975
$ws.Range("I34").Value = 37019.965
$ws.Range("J34").Value = 4627.273
$ws.Range("K34").Value = 37019.965
$ws.Range("L34").Value = 4627.273
$ws.Range("M34").Value = -36817.965
$ws.Range("N34").Value = -5031.273
$ws.Range("H110").Value = 29702
$ws.Range("J110").Value = 29702
$ws.Range("L110").Value = 29702
$ws.Range("N110").Value = -37882
$ws.Range("H134").Value = 953.4828
$ws.Range("I134").Value = 953.4828
$ws.Range("K134").Value = 2860.4484
$ws.Range("M134").Value = -325.4484000000002

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19989.2
$ws.Range("J46").Value = 19989.2
$ws.Range("L46").Value = 19989.2
$ws.Range("N46").Value = -20301.2
$ws.Range("H80").Value = 12530.1
$ws.Range("I80").Value = 2913.125
$ws.Range("J80").Value = 50998
$ws.Range("K80").Value = 2913.125
$ws.Range("L80").Value = 50998
$ws.Range("M80").Value = -1915.125
$ws.Range("N80").Value = -52994
$ws.Range("H83").Value = 12530.1
$ws.Range("I83").Value = 2913.125
$ws.Range("J83").Value = 50998
$ws.Range("K83").Value = 14565.625
$ws.Range("L83").Value = 254990
$ws.Range("M83").Value = -9573.625
$ws.Range("N83").Value = -264974

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 8000
$ws.Range("I18").Value = 8000
$ws.Range("K18").Value = 8000
$ws.Range("M18").Value = -7828
$ws.Range("H22").Value = 1035
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 1112.3684
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 1112.3684
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -1702.3684
$ws.Range("H27").Value = 1035
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 1112.3684
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 1112.3684
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -1326.3684
$ws.Range("H46").Value = 755
$ws.Range("I46").Value = 633.3333
$ws.Range("J46").Value = 876.6667
$ws.Range("K46").Value = 633.3333
$ws.Range("L46").Value = 876.6667
$ws.Range("M46").Value = -445.3333
$ws.Range("N46").Value = -1252.6667
$ws.Range("H68").Value = 2147.2222
$ws.Range("I68").Value = 1704.5454
$ws.Range("J68").Value = 2842.8572
$ws.Range("K68").Value = 1704.5454
$ws.Range("L68").Value = 2842.8572
$ws.Range("M68").Value = -955.5454
$ws.Range("N68").Value = -4340.8572
$ws.Range("H71").Value = 2147.2222
$ws.Range("I71").Value = 1704.5454
$ws.Range("J71").Value = 2842.8572
$ws.Range("K71").Value = 8522.726999999999
$ws.Range("L71").Value = 14214.286
$ws.Range("M71").Value = -4778.726999999999
$ws.Range("N71").Value = -21702.286
$ws.Range("H82").Value = 2365.5557
$ws.Range("I82").Value = 1466.6666
$ws.Range("J82").Value = 2815
$ws.Range("K82").Value = 1466.6666
$ws.Range("L82").Value = 2815
$ws.Range("M82").Value = -1105.6666
$ws.Range("N82").Value = -3537
$ws.Range("H85").Value = 2365.5557
$ws.Range("I85").Value = 1466.6666
$ws.Range("J85").Value = 2815
$ws.Range("K85").Value = 1466.6666
$ws.Range("L85").Value = 2815
$ws.Range("M85").Value = -218.6666
$ws.Range("N85").Value = -5311
$ws.Range("H136").Value = 346396.06
$ws.Range("I136").Value = 455496.62
$ws.Range("J136").Value = 3508.5715
$ws.Range("K136").Value = 1366489.86
$ws.Range("L136").Value = 10525.7145
$ws.Range("M136").Value = -1363939.86
$ws.Range("N136").Value = -15625.7145

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 200003900
$ws.Range("J62").Value = 333337000
$ws.Range("L62").Value = 333337000
$ws.Range("N62").Value = -333338248
$ws.Range("H65").Value = 200003900
$ws.Range("J65").Value = 333337000
$ws.Range("L65").Value = 1666685000
$ws.Range("N65").Value = -1666691240
$ws.Range("H74").Value = 5826.846
$ws.Range("J74").Value = 5826.846
$ws.Range("L74").Value = 5826.846
$ws.Range("N74").Value = -7698.846
$ws.Range("H77").Value = 5826.846
$ws.Range("J77").Value = 5826.846
$ws.Range("L77").Value = 17480.538
$ws.Range("N77").Value = -26840.538
$ws.Range("H132").Value = 2703.74
$ws.Range("I132").Value = 439.64444
$ws.Range("J132").Value = 23080.6
$ws.Range("K132").Value = 1318.93332
$ws.Range("L132").Value = 69241.79999999999
$ws.Range("M132").Value = 1211.06668
$ws.Range("N132").Value = -74301.79999999999
$ws.Range("H136").Value = 1663099.8
$ws.Range("I136").Value = 1832220.2
$ws.Range("K136").Value = 5496660.6
$ws.Range("M136").Value = -5494110.6
